# Update classification-report metrics for the first two blocks
# (target_col == previous_concussions) with the re-trained values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - label "0"
$ws.Range("B2").Value = 0.4
$ws.Range("C2").Value = 0.4
$ws.Range("D2").Value = 0.4000000000000001

# Row 3 - label "1"
$ws.Range("B3").Value = 0.5714285714285714
$ws.Range("C3").Value = 0.5714285714285714
$ws.Range("D3").Value = 0.5714285714285714

# Row 4 - label "accuracy"
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.5

# Row 5 - label "macro avg"
$ws.Range("B5").Value = 0.4857142857142857
$ws.Range("C5").Value = 0.4857142857142857
$ws.Range("D5").Value = 0.4857142857142858

# Row 6 - label "weighted avg"
$ws.Range("B6").Value = 0.5
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.5

# Row 7 - label "0"
$ws.Range("B7").Value = 0.5
$ws.Range("D7").Value = 0.5

# Row 8 - label "1"
$ws.Range("B8").Value = 0.6428571428571429
$ws.Range("C8").Value = 0.6428571428571429
$ws.Range("D8").Value = 0.6428571428571429

# Row 9 - label "accuracy"
$ws.Range("B9").Value = 0.5833333333333334
$ws.Range("C9").Value = 0.5833333333333334
$ws.Range("D9").Value = 0.5833333333333334
$ws.Range("E9").Value = 0.5833333333333334

# Row 10 - label "macro avg"
$ws.Range("B10").Value = 0.5714285714285714
$ws.Range("C10").Value = 0.5714285714285714
$ws.Range("D10").Value = 0.5714285714285714

# Row 11 - label "weighted avg"
$ws.Range("B11").Value = 0.5833333333333334
$ws.Range("C11").Value = 0.5833333333333334
$ws.Range("D11").Value = 0.5833333333333334
